$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current values for columns D, J, K, L, M, O, P across rows 2-21
# (the 2021 weekly data for this market/product got reshuffled onto different dates)
$snapshot = @{}
for ($r = 2; $r -le 21; $r++) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        J = $ws.Cells.Item($r, 10).Value2
        K = $ws.Cells.Item($r, 11).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
    }
}

# Map: new row -> source row (whose old values now populate the new row)
$rowMap = @{
    2 = 10
    3 = 17
    4 = 9
    5 = 7
    6 = 2
    7 = 20
    8 = 5
    9 = 4
    10 = 12
    11 = 16
    12 = 19
    13 = 11
    14 = 18
    15 = 3
    16 = 21
    17 = 8
    18 = 15
    19 = 13
    20 = 14
    21 = 6
}

foreach ($newRow in $rowMap.Keys) {
    $src = $snapshot[$rowMap[$newRow]]
    $ws.Cells.Item($newRow, 4).Value2 = $src.D
    $ws.Cells.Item($newRow, 10).Value2 = $src.J
    $ws.Cells.Item($newRow, 11).Value2 = $src.K
    $ws.Cells.Item($newRow, 12).Value2 = $src.L
    $ws.Cells.Item($newRow, 13).Value2 = $src.M
    $ws.Cells.Item($newRow, 15).Value2 = $src.O
    $ws.Cells.Item($newRow, 16).Value2 = $src.P
}
